$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build full data grid for A1:T23 and assign in one shot for speed/precision
$data = New-Object 'object[,]' 23,20
$data[0,0] = $null
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 2
$data[0,4] = 3
$data[0,5] = 4
$data[0,6] = 5
$data[0,7] = 6
$data[0,8] = 7
$data[0,9] = 8
$data[0,10] = 9
$data[0,11] = 10
$data[0,12] = 11
$data[0,13] = 12
$data[0,14] = 13
$data[0,15] = 14
$data[0,16] = 15
$data[0,17] = 16
$data[0,18] = 17
$data[0,19] = 18
$data[1,0] = 0
$data[1,1] = "HKL"
$data[1,2] = "[3, 2, 1]"
$data[1,3] = "[2, 2, 2]"
$data[1,4] = "[3, 1, 0]"
$data[1,5] = "[1, 1, 0]"
$data[1,6] = "[2, 2, 0]"
$data[1,7] = "[2, 0, 0]"
$data[1,8] = "[4, 0, 0]"
$data[1,9] = "[2, 1, 1]"
$data[1,10] = "1Pair-A"
$data[1,11] = "1Pair-B"
$data[1,12] = "2Pairs-A"
$data[1,13] = "2Pairs-B"
$data[1,14] = "3Pairs-A"
$data[1,15] = "3Pairs-B"
$data[1,16] = "3Pairs-C"
$data[1,17] = "4Pairs"
$data[1,18] = "5A4F"
$data[1,19] = "MaxUnique"
$data[2,0] = 1
$data[2,1] = "BT8Hex_2.5"
$data[2,2] = 0.9996053884200932
$data[2,3] = 0.9934522184426693
$data[2,4] = 1.001243156727804
$data[2,5] = 1.004079272695481
$data[2,6] = 1.004079272695481
$data[2,7] = 1.000056111302173
$data[2,8] = 1.000056111302173
$data[2,9] = 0.9978057941295956
$data[2,10] = 1.004079272695481
$data[2,11] = 0.9978057941295956
$data[2,12] = 0.9989309527158843
$data[2,13] = 0.9989309527158843
$data[2,14] = 0.9997016873865241
$data[2,15] = 1.00064705937575
$data[2,16] = 1.00064705937575
$data[2,17] = 1.001505112705683
$data[2,18] = 1.001505112705683
$data[2,19] = 0.9993736569529693
$data[3,0] = 2
$data[3,1] = "BT8Hex_5"
$data[3,2] = 0.9992476470130748
$data[3,3] = 0.9871908795186997
$data[3,4] = 1.002389873954126
$data[3,5] = 1.008044515889283
$data[3,6] = 1.008044515889283
$data[3,7] = 1.000037250976932
$data[3,8] = 1.000037250976932
$data[3,9] = 0.9957074297342631
$data[3,10] = 1.008044515889283
$data[3,11] = 0.9957074297342631
$data[3,12] = 0.9978723403555974
$data[3,13] = 0.9978723403555974
$data[3,14] = 0.9993781848884403
$data[3,15] = 1.001263065533493
$data[3,16] = 1.001263065533493
$data[3,17] = 1.00295842812244
$data[3,18] = 1.00295842812244
$data[3,19] = 0.9987695995143965
$data[4,0] = 3
$data[4,1] = "BT8Hex_10"
$data[4,2] = 0.9984809803450996
$data[4,3] = 0.9767254661281187
$data[4,4] = 1.004615044296095
$data[4,5] = 1.01424794905598
$data[4,6] = 1.01424794905598
$data[4,7] = 1.000781733067722
$data[4,8] = 1.000781733067722
$data[4,9] = 0.992156065828924
$data[4,10] = 1.01424794905598
$data[4,11] = 0.992156065828924
$data[4,12] = 0.9964688994483231
$data[4,13] = 0.9964688994483231
$data[4,14] = 0.9991842810642472
$data[4,15] = 1.002395249317542
$data[4,16] = 1.002395249317542
$data[4,17] = 1.005358424252152
$data[4,18] = 1.005358424252152
$data[4,19] = 0.99783453978699
$data[5,0] = 4
$data[5,1] = "BT8Hex_15"
$data[5,2] = 0.9976958765149293
$data[5,3] = 0.9665915891568944
$data[5,4] = 1.006790186687415
$data[5,5] = 1.02036364384464
$data[5,6] = 1.02036364384464
$data[5,7] = 1.001416785941297
$data[5,8] = 1.001416785941297
$data[5,9] = 0.9887093002183475
$data[5,10] = 1.02036364384464
$data[5,11] = 0.9887093002183475
$data[5,12] = 0.9950630430798222
$data[5,13] = 0.9950630430798222
$data[5,14] = 0.9989720909490197
$data[5,15] = 1.003496576668095
$data[5,16] = 1.003496576668095
$data[5,17] = 1.007713343462231
$data[5,18] = 1.007713343462231
$data[5,19] = 0.9969278970605872
$data[6,0] = 5
$data[6,1] = "Spiral2.5"
$data[6,2] = 0.9999508641556855
$data[6,3] = 1.000265861068155
$data[6,4] = 1.000052604629472
$data[6,5] = 0.9997256930910289
$data[6,6] = 0.9997256930910289
$data[6,7] = 1.00032785616558
$data[6,8] = 1.00032785616558
$data[6,9] = 1.0000482834371
$data[6,10] = 0.9997256930910289
$data[6,11] = 1.0000482834371
$data[6,12] = 1.00018806980134
$data[6,13] = 1.00018806980134
$data[6,14] = 1.000142914744051
$data[6,15] = 1.000033944231236
$data[6,16] = 1.000033944231236
$data[6,17] = 0.9999568814461843
$data[6,18] = 0.9999568814461843
$data[6,19] = 1.000061860424504
$data[7,0] = 6
$data[7,1] = "Spiral5"
$data[7,2] = 0.9999548459083727
$data[7,3] = 1.000635003917454
$data[7,4] = 0.9999645050297844
$data[7,5] = 0.9996359668358527
$data[7,6] = 0.9996359668358527
$data[7,7] = 1.000447202243416
$data[7,8] = 1.000447202243416
$data[7,9] = 1.000099052588053
$data[7,10] = 0.9996359668358527
$data[7,11] = 1.000099052588053
$data[7,12] = 1.000273127415735
$data[7,13] = 1.000273127415735
$data[7,14] = 1.000170253287085
$data[7,15] = 1.000060740555774
$data[7,16] = 1.000060740555774
$data[7,17] = 0.9999545471257938
$data[7,18] = 0.9999545471257938
$data[7,19] = 1.000122762753822
$data[8,0] = 7
$data[8,1] = "Spiral7.5"
$data[8,2] = 0.9999874852715748
$data[8,3] = 1.000417126450857
$data[8,4] = 0.9998723771757533
$data[8,5] = 1.000065026367057
$data[8,6] = 1.000065026367057
$data[8,7] = 1.000277662116666
$data[8,8] = 1.000277662116666
$data[8,9] = 0.999974651896293
$data[8,10] = 1.000065026367057
$data[8,11] = 0.999974651896293
$data[8,12] = 1.000126157006479
$data[8,13] = 1.000126157006479
$data[8,14] = 1.000041563729571
$data[8,15] = 1.000105780126672
$data[8,16] = 1.000105780126672
$data[8,17] = 1.000095591686768
$data[8,18] = 1.000095591686768
$data[8,19] = 1.0000990548797
$data[9,0] = 8
$data[9,1] = "Spiral10"
$data[9,2] = 0.9999120541462314
$data[9,3] = 1.001500274289551
$data[9,4] = 0.9998155960505235
$data[9,5] = 0.9994270165359904
$data[9,6] = 0.9994270165359904
$data[9,7] = 1.001016739602364
$data[9,8] = 1.001016739602364
$data[9,9] = 1.000152710179862
$data[9,10] = 0.9994270165359904
$data[9,11] = 1.000152710179862
$data[9,12] = 1.000584724891113
$data[9,13] = 1.000584724891113
$data[9,14] = 1.000328348610916
$data[9,15] = 1.000198822106072
$data[9,16] = 1.000198822106072
$data[9,17] = 1.000005870713552
$data[9,18] = 1.000005870713552
$data[9,19] = 1.000304065134087
$data[10,0] = 9
$data[10,1] = "Spiral15"
$data[10,2] = 1.000010299894898
$data[10,3] = 1.001218152282969
$data[10,4] = 0.9994022125434108
$data[10,5] = 1.000829849356118
$data[10,6] = 1.000829849356118
$data[10,7] = 1.000677738209495
$data[10,8] = 1.000677738209495
$data[10,9] = 0.9997432219009663
$data[10,10] = 1.000829849356118
$data[10,11] = 0.9997432219009663
$data[10,12] = 1.000210480055231
$data[10,13] = 1.000210480055231
$data[10,14] = 0.9999410575512907
$data[10,15] = 1.00041693648886
$data[10,16] = 1.00041693648886
$data[10,17] = 1.000520164705675
$data[10,18] = 1.000520164705675
$data[10,19] = 1.00031357903131
$data[11,0] = 10
$data[11,1] = "OffsetF45"
$data[11,2] = 0.9984295357935122
$data[11,3] = 0.8139450249381859
$data[11,4] = 1.025935258346493
$data[11,5] = 1.109063192977181
$data[11,6] = 1.109063192977181
$data[11,7] = 0.9494678571038061
$data[11,8] = 0.9494678571038061
$data[11,9] = 0.9527156541130148
$data[11,10] = 1.109063192977181
$data[11,11] = 0.9527156541130148
$data[11,12] = 0.9510917556084104
$data[11,13] = 0.9510917556084104
$data[11,14] = 0.9760395898544378
$data[11,15] = 1.003748901398001
$data[11,16] = 1.003748901398001
$data[11,17] = 1.030077474292796
$data[11,18] = 1.030077474292796
$data[11,19] = 0.9749260872120321
$data[12,0] = 11
$data[12,1] = "OffsetA45"
$data[12,2] = 0.9901037861370713
$data[12,3] = 1.059424617727472
$data[12,4] = 1.006722688232293
$data[12,5] = 0.9539315096628049
$data[12,6] = 0.9539315096628049
$data[12,7] = 1.073124679556563
$data[12,8] = 1.073124679556563
$data[12,9] = 1.005072393960247
$data[12,10] = 0.9539315096628049
$data[12,11] = 1.005072393960247
$data[12,12] = 1.039098536758405
$data[12,13] = 1.039098536758405
$data[12,14] = 1.028306587249701
$data[12,15] = 1.010709527726538
$data[12,16] = 1.010709527726538
$data[12,17] = 0.996515023210605
$data[12,18] = 0.996515023210605
$data[12,19] = 1.014729945879409
$data[13,0] = 12
$data[13,1] = "OffsetFTD"
$data[13,2] = 0.9572249217576906
$data[13,3] = 1.044483302036352
$data[13,4] = 1.072292437149501
$data[13,5] = 0.869042162737965
$data[13,6] = 0.869042162737965
$data[13,7] = 1.225090019647668
$data[13,8] = 1.225090019647668
$data[13,9] = 0.9932718048500814
$data[13,10] = 0.869042162737965
$data[13,11] = 0.9932718048500814
$data[13,12] = 1.109180912248874
$data[13,13] = 1.109180912248874
$data[13,14] = 1.096884753882417
$data[13,15] = 1.029134662411905
$data[13,16] = 1.029134662411905
$data[13,17] = 0.9891115374934196
$data[13,18] = 0.9891115374934196
$data[13,19] = 1.026900774696543
$data[14,0] = 13
$data[14,1] = "OffsetATD"
$data[14,2] = 0.9974445398336368
$data[14,3] = 0.9191395838574714
$data[14,4] = 1.021504866922262
$data[14,5] = 1.018881546054003
$data[14,6] = 1.018881546054003
$data[14,7] = 0.9747000332258716
$data[14,8] = 0.9747000332258716
$data[14,9] = 0.9887038079113074
$data[14,10] = 1.018881546054003
$data[14,11] = 0.9887038079113074
$data[14,12] = 0.9817019205685895
$data[14,13] = 0.9817019205685895
$data[14,14] = 0.9949695693531471
$data[14,15] = 0.9940951290637274
$data[14,16] = 0.9940951290637274
$data[14,17] = 1.000291733311296
$data[14,18] = 1.000291733311296
$data[14,19] = 0.9867290629674255
$data[15,0] = 14
$data[15,1] = "Holden2.5"
$data[15,2] = 0.9916903555962717
$data[15,3] = 0.8846791195182852
$data[15,4] = 1.024631648951119
$data[15,5] = 1.067517242934141
$data[15,6] = 1.067517242934141
$data[15,7] = 1.007898761581609
$data[15,8] = 1.007898761581609
$data[15,9] = 0.9611055859492399
$data[15,10] = 1.067517242934141
$data[15,11] = 0.9611055859492399
$data[15,12] = 0.9845021737654245
$data[15,13] = 0.9845021737654245
$data[15,14] = 0.9978786654939893
$data[15,15] = 1.01217386348833
$data[15,16] = 1.01217386348833
$data[15,17] = 1.026009708349783
$data[15,18] = 1.026009708349783
$data[15,19] = 0.9895871190884443
$data[16,0] = 15
$data[16,1] = "Holden5"
$data[16,2] = 0.9917352053653926
$data[16,3] = 0.9312323558947232
$data[16,4] = 1.020801670959478
$data[16,5] = 1.031745908548341
$data[16,6] = 1.031745908548341
$data[16,7] = 1.020685602925075
$data[16,8] = 1.020685602925075
$data[16,9] = 0.9757472300845478
$data[16,10] = 1.031745908548341
$data[16,11] = 0.9757472300845478
$data[16,12] = 0.9982164165048113
$data[16,13] = 0.9982164165048113
$data[16,14] = 1.005744834656367
$data[16,15] = 1.009392913852655
$data[16,16] = 1.009392913852655
$data[16,17] = 1.014981162526576
$data[16,18] = 1.014981162526576
$data[16,19] = 0.9953246622962597
$data[17,0] = 16
$data[17,1] = "Holden10"
$data[17,2] = 0.9919991496634896
$data[17,3] = 1.02344921829446
$data[17,4] = 1.012845847892839
$data[17,5] = 0.9615059119622862
$data[17,6] = 0.9615059119622862
$data[17,7] = 1.045207748605022
$data[17,8] = 1.045207748605022
$data[17,9] = 1.004742183351274
$data[17,10] = 0.9615059119622862
$data[17,11] = 1.004742183351274
$data[17,12] = 1.024974965978148
$data[17,13] = 1.024974965978148
$data[17,14] = 1.020931926616378
$data[17,15] = 1.003818614639528
$data[17,16] = 1.003818614639528
$data[17,17] = 0.9932404389702172
$data[17,18] = 0.9932404389702172
$data[17,19] = 1.006625009961562
$data[18,0] = 17
$data[18,1] = "Holden15"
$data[18,2] = 0.9905754699802765
$data[18,3] = 1.040944046825794
$data[18,4] = 1.013769105285833
$data[18,5] = 0.9454213127984917
$data[18,6] = 0.9454213127984917
$data[18,7] = 1.057136258940385
$data[18,8] = 1.057136258940385
$data[18,9] = 1.009424154582177
$data[18,10] = 0.9454213127984917
$data[18,11] = 1.009424154582177
$data[18,12] = 1.033280206761281
$data[18,13] = 1.033280206761281
$data[18,14] = 1.026776506269465
$data[18,15] = 1.003993908773684
$data[18,16] = 1.003993908773684
$data[18,17] = 0.9893507597798863
$data[18,18] = 0.9893507597798863
$data[18,19] = 1.009545058068826
$data[19,0] = 18
$data[19,1] = "HexGrid-90degTilt2.5degRes"
$data[19,2] = 0.999936542137644
$data[19,3] = 1.001484240236526
$data[19,4] = 1.000071926568565
$data[19,5] = 0.9983992234231776
$data[19,6] = 0.9983992234231776
$data[19,7] = 1.000640812260579
$data[19,8] = 1.000640812260579
$data[19,9] = 1.000572255497957
$data[19,10] = 0.9983992234231776
$data[19,11] = 1.000572255497957
$data[19,12] = 1.000606533879268
$data[19,13] = 1.000606533879268
$data[19,14] = 1.000428331442367
$data[19,15] = 0.999870763727238
$data[19,16] = 0.999870763727238
$data[19,17] = 0.999502878651223
$data[19,18] = 0.999502878651223
$data[19,19] = 1.000184166687408
$data[20,0] = 19
$data[20,1] = "HexGrid-90degTilt5degRes"
$data[20,2] = 0.9999919311780319
$data[20,3] = 0.9983497907190984
$data[20,4] = 1.000047358002037
$data[20,5] = 1.001644908188164
$data[20,6] = 1.001644908188164
$data[20,7] = 0.9996655245847991
$data[20,8] = 0.9996655245847991
$data[20,9] = 0.9993158025385912
$data[20,10] = 1.001644908188164
$data[20,11] = 0.9993158025385912
$data[20,12] = 0.9994906635616951
$data[20,13] = 0.9994906635616951
$data[20,14] = 0.9996762283751425
$data[20,15] = 1.000208745103851
$data[20,16] = 1.000208745103851
$data[20,17] = 1.000567785874929
$data[20,18] = 1.000567785874929
$data[20,19] = 0.9998358858684536
$data[21,0] = 20
$data[21,1] = "HexGrid-90degTilt10degRes"
$data[21,2] = 0.9997317362045347
$data[21,3] = 1.001653109060683
$data[21,4] = 1.000543946269374
$data[21,5] = 0.9972822023010329
$data[21,6] = 0.9972822023010329
$data[21,7] = 1.001447096620321
$data[21,8] = 1.001447096620321
$data[21,9] = 1.000795553994934
$data[21,10] = 0.9972822023010329
$data[21,11] = 1.000795553994934
$data[21,12] = 1.001121325307627
$data[21,13] = 1.001121325307627
$data[21,14] = 1.000928865628209
$data[21,15] = 0.9998416176387624
$data[21,16] = 0.9998416176387624
$data[21,17] = 0.9992017638043301
$data[21,18] = 0.9992017638043301
$data[21,19] = 1.000242274075146
$data[22,0] = 21
$data[22,1] = "HexGrid-90degTilt15degRes"
$data[22,2] = 0.9993609335869809
$data[22,3] = 1.007527058861167
$data[22,4] = 1.001545656278879
$data[22,5] = 0.9887443817820742
$data[22,6] = 0.9887443817820742
$data[22,7] = 1.004376702498299
$data[22,8] = 1.004376702498299
$data[22,9] = 1.003733753017434
$data[22,10] = 0.9887443817820742
$data[22,11] = 1.003733753017434
$data[22,12] = 1.004055227757867
$data[22,13] = 1.004055227757867
$data[22,14] = 1.003218703931537
$data[22,15] = 0.9989516124326027
$data[22,16] = 0.9989516124326027
$data[22,17] = 0.9963998047699705
$data[22,18] = 0.9963998047699705
$data[22,19] = 1.000881414337472

$ws.Range("A1:T23").Value2 = $data

# Re-apply header/index style (s="1") for the newly populated A-column rows (A16:A23)
# and ensure B1:T1 header-row style is preserved (already had style 1 originally and unaffected).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A23").PasteSpecial(-4122) | Out-Null

# Remove now-unused columns U:AD (old extra duplicate data) to shrink sheet dimension to A1:T23
$ws.Range("U1:AD23").EntireColumn.Delete() | Out-Null
